# Update values in Sheet1 per the diff (KNN algorithm result data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D7"  = -7.812
    "B10" = 5.777
    "B12" = 5.090000000000001
    "D15" = -8.175000000000001
    "B18" = 5.266
    "D20" = -7.562
    "D29" = -7.281999999999999
    "D30" = -7.178999999999999
    "D31" = -8.086
    "B37" = 8.73
    "D40" = -7.662999999999999
    "B55" = 4.621
    "B68" = 5.100999999999999
    "D68" = -6.778
    "D76" = -7.753
    "B77" = 5.880000000000001
    "B78" = 7.811
    "D87" = -8.113
    "D88" = -7.916000000000001
    "D96" = -7.267
    "D98" = -8.228
    "D101" = -7.886000000000001
    "D102" = -8.036
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
